$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.523.11'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '3.498.99'
$ws.Range('E3').Value = '  +0.42%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'602.95"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.57%  '
$ws.Range('D6').Value = "'194.02"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.16%  '
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('D8').Value = "'0.999"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = "'0.200"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.07%  '
$ws.Range('D10').Value = "'0.646"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.16%  '
$ws.Range('D11').Value = "'53.05"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.70%  '
$ws.Range('D12').Value = "'0.0000299"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.15%  '
$ws.Range('D13').Value = "'9.46"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.27%  '
$ws.Range('D14').Value = '4.057.75'
$ws.Range('E14').Value = '  +0.16%  '
$ws.Range('D15').Value = "'593.78"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.73%  '
$ws.Range('D16').Value = '69.719.73'
$ws.Range('E16').Value = '  +0.36%  '
$ws.Range('D17').Value = "'18.98"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.83%  '
$ws.Range('D18').Value = "'12.66"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.00%  '
$ws.Range('E19').Value = '  +2.37%  '
$ws.Range('D20').Value = '3.489.55'
$ws.Range('E20').Value = '  -1.04%  '
$ws.Range('D21').Value = "'0.983"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.31%  '
$ws.Range('D22').Value = "'18.23"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +7.04%  '
$ws.Range('D23').Value = "'5.27"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.00%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').Value = "'4.62"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.51%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = "'101.34"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.40%  '
$ws.Range('E26').Value = '  +4.42%  '
$ws.Range('D27').Value = "'10.78"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.96%  '
$ws.Range('D28').Value = "'9.47"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.62%  '
$ws.Range('D29').Value = "'33.06"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.64%  '
$ws.Range('E30').Value = '  +7.84%  '
$ws.Range('D31').Value = "'7.00"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.69%  '
$ws.Range('D32').Value = "'12.33"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.26%  '
$ws.Range('E33').Value = '  +0.14%  '
$ws.Range('D34').Value = "'63.06"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.12%  '
$ws.Range('D35').Value = '3.735.57'
$ws.Range('E35').Value = '  +3.35%  '
$ws.Range('E36').Value = '  -0.36%  '
$ws.Range('D37').Value = '0.0₃0812'
$ws.Range('E37').Value = '  +5.89%  '
$ws.Range('E38').Value = '  +0.11%  '
$ws.Range('D39').Value = "'3.63"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.77%  '
$ws.Range('D40').Value = "'0.389"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.63%  '
$ws.Range('D41').Value = "'36.19"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.63%  '
$ws.Range('D42').Value = "'490.88"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.59%  '
$ws.Range('D43').Value = "'0.134"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.20%  '
$ws.Range('D44').Value = "'0.0450"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.68%  '
$ws.Range('D45').Value = "'0.139"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.73%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').Value = "'3.29"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('B47').Value = 'ThetaToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D47').Value = "'2.79"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.04%  '
$ws.Range('D48').Value = "'1.01"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.33%  '
$ws.Range('E49').Value = '  -3.53%  '
$ws.Range('D50').Value = "'0.000243"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.29%  '
$ws.Range('E51').Value = '  +10.28%  '
